# Rename the existing sheet "Sheet1" to "WisGate Edge "
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "WisGate Edge "

# Leave the selection on the original sheet where the author last left it
$ws1.Range("B33").Select()

# Add a new worksheet after the existing one for the WisGate Developer models
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws2 = $wb.Worksheets.Add($null, $lastSheet)
$ws2.Name = "WisGate Developer "

# Fill in the Description column first (column B, rows 2-6)
$ws2.Range("B2").Value = "WisGate Developer D4 / D4+ / D4P"
$ws2.Range("B3").Value = "WisGate Developer D0 / D0+"
$ws2.Range("B4").Value = "WisGate Developer D3 / D3+"
$ws2.Range("B5").Value = "WisGate Developer D4H"
$ws2.Range("B6").Value = "WisGate Developer Base"

# Fill in the model name column next (column A, rows 2-6)
$ws2.Range("A2").Value = "RAK7244 / RAK7244C / RAK7244P"
$ws2.Range("A3").Value = "RAK7246 / RAK7246G"
$ws2.Range("A4").Value = "RAK7243 / RAK7243C"
$ws2.Range("A5").Value = "RAK7248"
$ws2.Range("A6").Value = "RAK7271 / RAK7371"

# Fill in the header row last
$ws2.Range("A1").Value = "model name"
$ws2.Range("B1").Value = "Description"

# Apply a thin border around every cell of the table
$rng = $ws2.Range("A1:B6")
for ($i = 1; $i -le 4; $i++) {
    $rng.Borders.Item($i).LineStyle = 1
    $rng.Borders.Item($i).Weight = 2
}

# Auto-fit the columns to the new content
$rng.EntireColumn.AutoFit()

# Make the new sheet the active tab, with the cursor where the author left it
$ws2.Activate()
$ws2.Range("H13").Select()
